# Overworld 1 gameplay revision: add a new "Topographic Features" section
# to the language table, inserted as 9 new rows (30-38) before the
# existing "Season" section, shifting all subsequent rows down by 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows at row 30 (pushes old rows 30-198 down to 39-207)
$ws.Rows("30:38").Insert()

# Populate the new rows with the Topographic Features strings, in the
# same entry order the original author used (this controls the order
# new entries are appended to the shared string table).
$ws.Range("A30").Value = "topographyFeatures"
$ws.Range("A33").Value = "topography_forest"
$ws.Range("A34").Value = "topography_jungle"
$ws.Range("A37").Value = "topography_lake"
$ws.Range("A38").Value = "topography_ocean"
$ws.Range("B33").Value = "Forest"
$ws.Range("B34").Value = "Jungle"
$ws.Range("B37").Value = "Lake"
$ws.Range("B38").Value = "Ocean"
$ws.Range("B30").Value = "Topographic Features"
$ws.Range("A36").Value = "topography_river"
$ws.Range("B36").Value = "River"
$ws.Range("A31").Value = "topography_mountain"
$ws.Range("A32").Value = "topography_hill"
$ws.Range("B31").Value = "Mountain"
$ws.Range("B32").Value = "Hill"
$ws.Range("B35").Value = "Vegetation"
$ws.Range("A35").Value = "topography_vegetation"

# Update the view state to match the new selection/scroll position
$ws.Range("A16").Select()
$ws.Range("A34").Select()
